# Updates cryptos list values/percentages per the Wed Apr 24 06:44:31 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text reads as a plain decimal number need NumberFormat forced to
# Text first, otherwise Excel auto-converts the literal into a numeric value (losing
# trailing zeros / exact decimal text, e.g. "494.40" -> 494.4, "0.0000272" -> 2.72E-05).
# ClearFormats() afterwards restores the original (default, unstyled) cell format so
# only the cell content changes.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

$ws.Range("D2").Value = "66.692.94"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "3.254.44"
$ws.Range("E3").Value = "  +2.22%  "
Set-TextValue "D4" "0.998"
$ws.Range("E4").Value = "  -0.22%  "
Set-TextValue "D5" "606.31"
$ws.Range("E5").Value = "  +0.80%  "
Set-TextValue "D6" "157.99"
$ws.Range("E6").Value = "  +1.91%  "
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.255.62"
$ws.Range("E8").Value = "  +2.36%  "
Set-TextValue "D9" "0.549"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  +1.74%  "
Set-TextValue "D11" "5.95"
$ws.Range("E11").Value = "  +6.58%  "
Set-TextValue "D12" "0.507"
$ws.Range("E12").Value = "  -1.38%  "
Set-TextValue "D13" "0.0000272"
$ws.Range("E13").Value = "  +1.46%  "
Set-TextValue "D14" "39.34"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "3.780.09"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").Value = "66.662.26"
$ws.Range("E16").Value = "  +0.22%  "
Set-TextValue "D17" "7.42"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "3.245.99"
$ws.Range("E18").Value = "  +1.99%  "
Set-TextValue "D19" "0.114"
$ws.Range("E19").Value = "  +1.27%  "
Set-TextValue "D20" "509.22"
$ws.Range("E20").Value = "  -0.79%  "
Set-TextValue "D21" "15.43"
$ws.Range("E21").Value = "  -0.22%  "
Set-TextValue "D22" "0.752"
$ws.Range("E22").Value = "  +2.49%  "
Set-TextValue "D23" "8.12"
$ws.Range("E23").Value = "  -0.30%  "
Set-TextValue "D24" "14.87"
$ws.Range("E24").Value = "  -0.26%  "
Set-TextValue "D25" "86.69"
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("E26").Value = "  +0.27%  "
Set-TextValue "D27" "0.143"
$ws.Range("E27").Value = "  +60.36%  "
Set-TextValue "D28" "3.02"
$ws.Range("E28").Value = "  +0.54%  "
Set-TextValue "D29" "9.05"
$ws.Range("E29").Value = "  -2.02%  "
Set-TextValue "D30" "2.39"
$ws.Range("E30").Value = "  -0.30%  "
Set-TextValue "D31" "2.87"
$ws.Range("E31").Value = "  -8.20%  "
Set-TextValue "D32" "6.85"
$ws.Range("E32").Value = "  -2.36%  "
Set-TextValue "D33" "28.12"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -4.30%  "
Set-TextValue "D36" "6.44"
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0791"
$ws.Range("E37").Value = "  +15.82%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D38" "55.47"
$ws.Range("E38").Value = "  +1.15%  "
Set-TextValue "D39" "3.32"
$ws.Range("E39").Value = "  +18.92%  "
Set-TextValue "D40" "494.40"
$ws.Range("E40").Value = "  -4.03%  "
Set-TextValue "D41" "0.0427"
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("E42").Value = "  +0.96%  "
Set-TextValue "D43" "8.84"
$ws.Range("E43").Value = "  -0.28%  "
Set-TextValue "D44" "0.294"
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("D45").Value = "2.963.43"
$ws.Range("E45").Value = "  +4.00%  "
Set-TextValue "D46" "2.48"
$ws.Range("E46").Value = "  +1.07%  "
Set-TextValue "D47" "28.71"
$ws.Range("E47").Value = "  +1.87%  "
Set-TextValue "D48" "2.50"
$ws.Range("E48").Value = "  +4.42%  "
$ws.Range("E49").Value = "  +2.66%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
Set-TextValue "D51" "2.54"
$ws.Range("E51").Value = "  -1.63%  "
